$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 28 (shifts old rows 28-218 down to 31-221)
$ws.Range("A28:A30").EntireRow.Insert()

# G27: add a note about using fewer features
$ws.Range("G27").Value = "Used less features to stop overfitting"

# Row 28: new experiment result
$ws.Range("A28").Value = "Title_Mr, ""Sex"", ""Title_Mrs"", ""Pclass_3"", ""Title_Miss"", ""Cabin_NA"", ""Fare_0"""
$ws.Range("B28").Value = "{'criterion': 'entropy', 'max_depth': 30, 'max_features': 'sqrt', 'min_samples_leaf': 10, 'min_samples_split': 10, 'n_estimators': 10}"
$ws.Range("C28").Value = 82.019999999999897
$ws.Range("D28").Value = 80.45
$ws.Range("E28").Value = 80.360799001248395
$ws.Range("F28").Value = 78.947000000000003

# Row 29: new experiment result + note
$ws.Range("A29").Value = "Title_Mr, ""Sex"", ""Title_Mrs"", ""Pclass_3"", ""Title_Miss"""
$ws.Range("B29").Value = "{'criterion': 'gini', 'max_depth': 3, 'max_features': 'sqrt', 'min_samples_leaf': 10, 'min_samples_split': 30, 'n_estimators': 10}"
$ws.Range("C29").Value = 79.069999999999894
$ws.Range("D29").Value = 78.209999999999894
$ws.Range("E29").Value = 78.1111111111111
$ws.Range("G29").Value = "Not going to test this model"

# Row 30: new experiment result
$ws.Range("A30").Value = "Title_Mr, ""Sex"", ""Title_Mrs"", ""Pclass_3"", ""Title_Miss"", ""Cabin_NA"", ""Fare_0"", ""Age_3"", ""Age_2"", ""Embarked_C"", ""Age_4"", ""LargeFamily"", ""Pclass_1"", ""Pclass_2"", ""Embarked_S"""
$ws.Range("B30").Value = "{'criterion': 'entropy', 'max_depth': 100, 'max_features': 'sqrt', 'min_samples_leaf': 3, 'min_samples_split': 3, 'n_estimators': 100}"
$ws.Range("C30").Value = 84.829999999999899
$ws.Range("D30").Value = 83.799999999999898
$ws.Range("E30").Value = 83.277153558052404
$ws.Range("F30").Value = 77.510999999999996

# Small style tweak further down: G161:G163 pick up the "no-wrap" style used by G2:G160
$ws.Range("G161:G163").Style = $ws.Range("G160").Style

# Update the sheet view (scrolled position / active selection)
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("G30").Select()
